# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet named "2022-Q4" right before the existing
#    "2022-Q3" sheet (so it becomes the 2nd sheet, pushing all the
#    quarterly sheets down by one position) and populate it with the
#    fund-holding detail table for that quarter.
# 2) Insert a new row in the "总计" (summary) sheet for "2022-Q4" right
#    after the header row, shifting the previously existing summary
#    rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert the 2022-Q4 summary row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Reapply the same formatting the other index cells in column A use
# (bold, centered, bordered) by copying the format from the row below.
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 0.27

# ---------------------------------------------------------------------
# Part 2: brand new "2022-Q4" fund-holding detail sheet
# ---------------------------------------------------------------------
$target = $wb.Worksheets.Item("2022-Q3")
$formatSrc = $wb.Worksheets.Item("2021-Q1")

$new = $wb.Worksheets.Add($target)
$new.Name = "2022-Q4"

# Outline / page setup metadata to match the other sheets.
$new.Outline.SummaryRow = 1
$new.Outline.SummaryColumn = 1
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

# Copy over the header-row formatting (bold/centered/bordered) and the
# generic data-row formatting (column A bold/centered/bordered, the
# rest unstyled), tiling the one-row pattern down across all 9 rows.
$formatSrc.Range("A1:H1").Copy()
$new.Range("A1:H1").PasteSpecial(-4122)
$formatSrc.Range("A2:H2").Copy()
$new.Range("A2:H9").PasteSpecial(-4122)

# ---- header row ----
$new.Cells.Item(1, 2).Value = "基金代码"
$new.Cells.Item(1, 3).Value = "基金名称"
$new.Cells.Item(1, 4).Value = "基金规模"
$new.Cells.Item(1, 5).Value = "股票总仓位"
$new.Cells.Item(1, 6).Value = "仓位占比"
$new.Cells.Item(1, 7).Value = "持有市值(亿元)"
$new.Cells.Item(1, 8).Value = "仓位排名"

# ---- data rows ----
# Columns B (fund code) and D/E/F/G (numeric-looking text figures) must
# stay text, so their number format is forced to Text ("@") before the
# value is assigned - otherwise COM would silently coerce strings like
# "001167" or "0.1385" into numbers and lose the original formatting.
$rows = @(
    @{r=2; a=0; b="001167"; c="金鹰科技创新股票";            d="3.17"; e="91.02"; f="4.37"; g="0.1385"; h=9},
    @{r=3; a=1; b="210002"; c="金鹰红利价值混合A";           d="0.96"; e="77.22"; f="5.94"; g="0.0570"; h=3},
    @{r=4; a=2; b="501073"; c="华安智联混合（LOF）A";        d="3.27"; e="36.94"; f="1.14"; g="0.0373"; h=8},
    @{r=5; a=3; b="016563"; c="金鹰红利价值混合C";           d="0.52"; e="77.22"; f="5.94"; g="0.0309"; h=3},
    @{r=6; a=4; b="011150"; c="创金合信ESG责任投资股票C";     d="0.12"; e="90.04"; f="2.74"; g="0.0033"; h=2},
    @{r=7; a=5; b="011149"; c="创金合信ESG责任投资股票A";     d="0.10"; e="90.04"; f="2.74"; g="0.0027"; h=2},
    @{r=8; a=6; b="001978"; c="泰信互联网+主题灵活配置混合"; d="0.06"; e="77.79"; f="2.63"; g="0.0016"; h=3},
    @{r=9; a=7; b="016071"; c="华安智联混合（LOF）C";        d="0.00"; e="36.94"; f="1.14"; g="__NUM0__"; h=8}
)

foreach ($row in $rows) {
    $new.Cells.Item($row.r, 1).Value = $row.a

    $new.Cells.Item($row.r, 2).NumberFormat = "@"
    $new.Cells.Item($row.r, 2).Value = $row.b

    $new.Cells.Item($row.r, 3).Value = $row.c

    $new.Cells.Item($row.r, 4).NumberFormat = "@"
    $new.Cells.Item($row.r, 4).Value = $row.d

    $new.Cells.Item($row.r, 5).NumberFormat = "@"
    $new.Cells.Item($row.r, 5).Value = $row.e

    $new.Cells.Item($row.r, 6).NumberFormat = "@"
    $new.Cells.Item($row.r, 6).Value = $row.f

    if ($row.g -eq "__NUM0__") {
        $new.Cells.Item($row.r, 7).Value = 0
    } else {
        $new.Cells.Item($row.r, 7).NumberFormat = "@"
        $new.Cells.Item($row.r, 7).Value = $row.g
    }

    $new.Cells.Item($row.r, 8).Value = $row.h
}

# Keep the originally-active/selected sheet ("总计") active, matching
# the rest of the unchanged workbook view state.
$summary.Activate()
